$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("begroting")

# Add the three new label rows under the existing table
$ws.Range("A21").Value = "hours"
$ws.Range("B21").Formula = "=24480/135"

$ws.Range("A22").Value = "booked"
$ws.Range("B22").Formula = "=8+34"

$ws.Range("A23").Value = "remaining"
$ws.Range("B23").Formula = "=B21-B22"
$ws.Range("C23").Formula = "=(B21-B22)/B21*100"
$ws.Range("C23").NumberFormat = "0.0"

# Make "begroting" the active sheet with the new selection
$ws.Activate()
$ws.Range("C24").Select()

# Restore the "api calls" sheet's prior (non-active) selection
$ws2 = $wb.Worksheets.Item("api calls")
$ws2.Range("B12").Select()
